$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "horizontal center" alignment that was applied to the
# non-header rows of column A (A2:A6, A10:A17) -- they keep their
# vertical centering but drop the explicit horizontal centering.
$ws.Range("A2:A6").HorizontalAlignment = 1
$ws.Range("A10:A17").HorizontalAlignment = 1

# Unmerge the two merged blocks in column A.
$ws.Range("A7:A9").UnMerge()
$ws.Range("A18:A19").UnMerge()

# The cells that used to be covered by the merge (and were blank) now
# get a single space as their literal value.
$ws.Range("A8").Value = " "
$ws.Range("A9").Value = " "
$ws.Range("A19").Value = " "

# Update the active selection / window state to match the saved view.
$ws.Range("H20").Select()
